$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1450151057401813
$ws.Range("C2").Value = 0.6676737160120846
$ws.Range("J2").Value = 0.01510574018126888
$ws.Range("P2").Value = 0.1148036253776435
$ws.Range("S2").Value = 0.05740181268882175
$ws.Range("B3").Value = 0.004385964912280702
$ws.Range("C3").Value = 0.03070175438596491
$ws.Range("J3").Value = 0.05263157894736842
$ws.Range("P3").Value = 0.7456140350877193
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.03125
$ws.Range("P4").Value = 0.671875
$ws.Range("S4").Value = 0.296875
$ws.Range("J5").Value = 1
$ws.Range("B6").Value = 0.08520179372197309
$ws.Range("D6").Value = 0.0179372197309417
$ws.Range("F6").Value = 0.08968609865470852
$ws.Range("J6").Value = 0.2914798206278027
$ws.Range("O6").Value = 0.01345291479820628
$ws.Range("Q6").Value = 0.1165919282511211
$ws.Range("R6").Value = 0.07174887892376682
$ws.Range("S6").Value = 0.3139013452914798
$ws.Range("B7").Value = 0.1149425287356322
$ws.Range("D7").Value = 0.01149425287356322
$ws.Range("F7").Value = 0.02298850574712644
$ws.Range("J7").Value = 0.1206896551724138
$ws.Range("O7").Value = 0.04022988505747126
$ws.Range("Q7").Value = 0.2011494252873563
$ws.Range("R7").Value = 0.08620689655172414
$ws.Range("S7").Value = 0.4022988505747127
$ws.Range("B8").Value = 0.1120879120879121
$ws.Range("D8").Value = 0.02637362637362637
$ws.Range("F8").Value = 0.04615384615384616
$ws.Range("J8").Value = 0.1208791208791209
$ws.Range("O8").Value = 0.01318681318681319
$ws.Range("Q8").Value = 0.2043956043956044
$ws.Range("R8").Value = 0.0989010989010989
$ws.Range("S8").Value = 0.378021978021978
$ws.Range("B9").Value = 0.1370558375634518
$ws.Range("D9").Value = 0.04060913705583756
$ws.Range("F9").Value = 0.05076142131979695
$ws.Range("J9").Value = 0.06598984771573604
$ws.Range("O9").Value = 0.02030456852791878
$ws.Range("Q9").Value = 0.1065989847715736
$ws.Range("R9").Value = 0.116751269035533
$ws.Range("S9").Value = 0.4619289340101523
$ws.Range("B10").Value = 0.1167146974063401
$ws.Range("D10").Value = 0.02737752161383285
$ws.Range("E10").Value = 0.0007204610951008645
$ws.Range("F10").Value = 0.06268011527377522
$ws.Range("J10").Value = 0.1376080691642651
$ws.Range("O10").Value = 0.01729106628242075
$ws.Range("Q10").Value = 0.2219020172910663
$ws.Range("R10").Value = 0.08717579250720461
$ws.Range("S10").Value = 0.3285302593659942
$ws.Range("G11").Value = 0.1529850746268657
$ws.Range("J11").Value = 0.0708955223880597
$ws.Range("K11").Value = 0.2052238805970149
$ws.Range("L11").Value = 0.5485074626865671
$ws.Range("S11").Value = 0.02238805970149254
$ws.Range("G12").Value = 0.7516339869281046
$ws.Range("J12").Value = 0.196078431372549
$ws.Range("L12").Value = 0.0261437908496732
$ws.Range("S12").Value = 0.0261437908496732
$ws.Range("F15").Value = 0.007751937984496124
$ws.Range("H15").Value = 0.1395348837209302
$ws.Range("I15").Value = 0.06589147286821706
$ws.Range("J15").Value = 0.3837209302325582
$ws.Range("K15").Value = 0.06976744186046512
$ws.Range("M15").Value = 0.007751937984496124
$ws.Range("O15").Value = 0.09302325581395349
$ws.Range("S15").Value = 0.2325581395348837
$ws.Range("F16").Value = 0.02880658436213992
$ws.Range("H16").Value = 0.1810699588477366
$ws.Range("I16").Value = 0.06172839506172839
$ws.Range("J16").Value = 0.4074074074074074
$ws.Range("K16").Value = 0.102880658436214
$ws.Range("M16").Value = 0.03292181069958848
$ws.Range("O16").Value = 0.08230452674897119
$ws.Range("S16").Value = 0.102880658436214
$ws.Range("F17").Value = 0.02685950413223141
$ws.Range("H17").Value = 0.1694214876033058
$ws.Range("I17").Value = 0.09504132231404959
$ws.Range("J17").Value = 0.4421487603305785
$ws.Range("K17").Value = 0.08057851239669421
$ws.Range("M17").Value = 0.01652892561983471
$ws.Range("O17").Value = 0.05991735537190083
$ws.Range("S17").Value = 0.109504132231405
$ws.Range("F18").Value = 0.02727272727272727
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.06818181818181818
$ws.Range("J18").Value = 0.4136363636363636
$ws.Range("K18").Value = 0.07272727272727272
$ws.Range("M18").Value = 0.01818181818181818
$ws.Range("O18").Value = 0.08636363636363636
$ws.Range("S18").Value = 0.1136363636363636
$ws.Range("F19").Value = 0.02247191011235955
$ws.Range("H19").Value = 0.2014446227929374
$ws.Range("I19").Value = 0.08426966292134831
$ws.Range("J19").Value = 0.3844301765650081
$ws.Range("K19").Value = 0.08828250401284109
$ws.Range("M19").Value = 0.01043338683788122
$ws.Range("O19").Value = 0.07784911717495988
$ws.Range("S19").Value = 0.1308186195826645
